{"js": "// Apply the LOT2025.docx content updates:\n// 1) Activation date 2018 -> 2025\n// 2) Objectives paragraph (PT) rewritten\n// 3) Objectives paragraph (EN, italic) - was empty, now has translated text\n// 4) \"Programa resumido\" (PT) rewritten\n// 5) \"Programa resumido\" (EN, italic) rewritten\n// 6) \"Programa\" (PT) rewritten\n// 7) \"Programa\" (EN, italic) rewritten\n// 8) Grading formula: \"(P1 + P2 )/2\" -> \"(P1 + P2)/2.\"\n// 9) Bibliografia paragraph rewritten\n\nconst replacements = [\n  {\n    old: \"Ativa\u00e7\u00e3o: 01/01/2018\",\n    new: \"Ativa\u00e7\u00e3o: 01/01/2025\"\n  },\n  {\n    old: \"Capacitar o aluno para identificar e aplicar os conceitos de Reatores Bioqu\u00edmicos em diferentes bioprocessos (fermentativos e enzim\u00e1ticos). Especificamente, capacitar o aluno para definir os tipos de biorreatores; para definir as diferentes formas de condu\u00e7\u00e3o de um processo fermentativo empregando biorreatores, para realizar o equacionamento matem\u00e1tico do crescimento microbiano e da forma\u00e7\u00e3o de produtos de interesse de um bioprocesso empregando biorreatores e para aplica\u00e7\u00e3o dos conceitos gerais dos reatores enzim\u00e1ticos operados em diferentes fases.\",\n    new: \"Desenvolver nos discentes as compet\u00eancias e habilidades necess\u00e1rias para a aplica\u00e7\u00e3o de conhecimentos cient\u00edficos, tecnol\u00f3gicos e de engenharia na concep\u00e7\u00e3o, projeto, instala\u00e7\u00e3o, otimiza\u00e7\u00e3o, supervis\u00e3o e avalia\u00e7\u00e3o cr\u00edtica da opera\u00e7\u00e3o de bioprocessos, com \u00eanfase em: 1) Tipos de biorreatores; 2) Formas de opera\u00e7\u00e3o dos biorreatores e 3) An\u00e1lise de biorreatores\"\n  },\n  {\n    old: \"Introdu\u00e7\u00e3o a biorreatores; processo descont\u00ednuo; processo cont\u00ednuo; processo descont\u00ednuo alimentado e reatores enzim\u00e1ticos.\",\n    new: \"Ser\u00e3o apresentados os principais tipos de biorreatores associados as suas aplica\u00e7\u00f5es para diferentes bioprocessos; definidas as formas de opera\u00e7\u00e3o do biorreator e analisadas as diferentes formas de opera\u00e7\u00e3o de biorreatores com base nos balan\u00e7os materiais dos componentes do sistema.\"\n  },\n  {\n    old: \"Introduction to bioreactors; batch process; continuous process; fed batch process and enzymatic reactors.\",\n    new: \"The main types of bioreactors associated with their applications for different bioprocesses will be presented. Also, the main forms of bioreactor operation will be defined and further analyzed based on the material balances of the system components.\"\n  },\n  {\n    old: \"1. Introdu\u00e7\u00e3o a biorreatores: apresenta\u00e7\u00e3o e classifica\u00e7\u00e3o de reatores bioqu\u00edmicos; intera\u00e7\u00e3o microorganismos/meios para estabelecimento de condi\u00e7\u00f5es para c\u00e1lculos de biorreatores.2. Processo descont\u00ednuo: caracter\u00edsticas gerais do processo cont\u00ednuo; balan\u00e7o material para c\u00e9lula, substrato e produto, em um \u00fanico est\u00e1gio com e sem reciclo de c\u00e9lulas; aplica\u00e7\u00e3o do processo cont\u00ednuo (exemplos).3. Processo cont\u00ednuo: caracter\u00edsticas gerais do processo cont\u00ednuo; balan\u00e7o material para c\u00e9lula, substrato e produto, em um \u00fanico est\u00e1gio com e sem reciclo de c\u00e9lulas; aplica\u00e7\u00e3o do processo cont\u00ednuo (exemplos).4. Processo descont\u00ednuo alimentado: caracter\u00edsticas gerais do processo descont\u00ednuo alimentado; balan\u00e7o material para c\u00e9lula e substrato, com volume vari\u00e1vel, empregando vaz\u00e3o constante de alimenta\u00e7\u00e3o; considera\u00e7\u00f5es sobre forma\u00e7\u00e3o de produtos no processo descont\u00ednuo alimentado; aplica\u00e7\u00e3o do processo descont\u00ednuo alimentado (exemplos).5. Reatores enzim\u00e1ticos: caracter\u00edsticas gerais dos reatores enzim\u00e1ticos; aplica\u00e7\u00e3o de processos enzim\u00e1ticos (exemplos).\",\n    new: \"1) Defini\u00e7\u00e3o e classifica\u00e7\u00e3o de biorreatores; 2) Biorreatores para cultivos submersos (agitados mecanicamente, agitados pneumaticamente, leito fixo, leito fluidizado, outros tipos; 3) Biorreatores para cultivos em estado s\u00f3lido (est\u00e1ticos e agitados); 4) Principais formas de opera\u00e7\u00e3o de biorreatores (descont\u00ednua, cont\u00ednua e descont\u00ednua-alimentada) e 4) An\u00e1lise de biorreatores (balan\u00e7os materiais para c\u00e9lulas, substrato-limitante e produtos metab\u00f3licos) nas diferentes formas de opera\u00e7\u00e3o.\"\n  },\n  {\n    old: \"1.Introduction to bioreactors: Description and classification of biochemical reactors; bioreactors configuration; microorganisms/medium interaction; bioreactor operation modes.2.Batch process: general characteristics of the batch process; material balance for cell, substrate and product; application of batch process (examples).3.Continuous process: general characteristics, operation modes (single or multiple stages with or without cells recycle); material balance for cell and substrate; formation of products in the continuous systems, application of continuous process (examples).4.Fed batch process: general characteristics; material balance for cell and substrate (equations for fed-batch operation at variable and fixed volume); considerations about formation of products in the fed batch process; application of fed batch process (examples).5.Enzymatic reactors: general characteristics of the enzymatic reactors; application of enzymatic processes (examples).\",\n    new: \"1) Definition and classification of bioreactors; 2) Bioreactors for submerged cultures (mechanically agitated, pneumatically agitated, fixed bed, fluidized bed, other types; 3) Bioreactors for solid state cultures (static and agitated); 4) Main forms of bioreactor operation (batch, continuous and fed-batch) and 4) Analysis of bioreactors (material balances for cells, limiting substrate and metabolic products) in different forms of operation.\"\n  },\n  {\n    old: \"Os alunos ser\u00e3o avaliados formalmente por duas provas te\u00f3ricas. A pondera\u00e7\u00e3o das notas ser\u00e1 de 50% para cada avalia\u00e7\u00e3o, ou seja: M\u00e9dia do per\u00edodo letivo normal = (P1 + P2 )/2\",\n    new: \"Os alunos ser\u00e3o avaliados formalmente por duas provas te\u00f3ricas. A pondera\u00e7\u00e3o das notas ser\u00e1 de 50% para cada avalia\u00e7\u00e3o, ou seja: M\u00e9dia do per\u00edodo letivo normal = (P1 + P2)/2.\"\n  },\n  {\n    old: \"1)  Aiba, S., Humphrey, A.E., Millis, N.F. Biochemical Engineering - 2\u00aa Edi\u00e7\u00e3o- 1973.2) Asenjo A., Merchuk, J.C. Bioreactor System Design-1995.3) Stanbury, D. and Whitaker, A. Principles af Fermentation Technology-1986.4) Lima, U.A., Aquarone, E., Borzani, W. Biotecnologia Industrial. Fundamentos Vol. 1, Engenharia Bioqu\u00edmica Vol.2, Processos Fermentativos Vol.3. Ed.  Edgard Blucher, S\u00e3o Paulo, 2001.\",\n    new: \"ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914PESSOA JR, A; VITOLO, M; LONG, P.F.(editors).  Pharmaceutical Biotechnology: A Focus on Industrial Application. CRC Press, 1st Edition.2022.\"\n  }\n];\n\nconst enObjectives = \"Develop in students the skills and abilities necessary to apply scientific, technological and engineering knowledge in the conception, design, installation, optimization, supervision and critical evaluation of the operation of bioprocesses, with an emphasis on: 1) Types of bioreactors; 2) Bioreactor operating modes and 3) Bioreactor analysis.\";\n\n// The English objectives paragraph is an (initially empty) italic run two\n// paragraphs after the \"Objetivos\" Heading 2 (heading, PT text, EN text).\n// Locate it this way - independent of the PT text replacement above/below -\n// since there is no existing text in it to search for.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length - 2; i++) {\n  if (paragraphs.items[i].style === \"Heading 2\" && paragraphs.items[i].text === \"Objetivos\") {\n    targetParagraph = paragraphs.items[i + 2];\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not locate empty English objectives paragraph\");\n}\ntargetParagraph.insertText(enObjectives, Word.InsertLocation.replace);\nawait context.sync();\n\nfor (const { old, new: replacement } of replacements) {\n  const results = context.document.body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + old.substring(0, 60));\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert the (until now empty) English objectives paragraph text.\n#    Located relative to the 'Objetivos' Heading 2 paragraph so this step\n#    does not depend on the Portuguese objectives text (which we rewrite below).\n$paragraphs = $d.Paragraphs\n$enObjectivesParagraph = $null\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $cur = $paragraphs.Item($i)\n    if ($cur.Style.NameLocal -eq \"Heading 2\" -and $cur.Range.Text.TrimEnd([char]13) -eq \"Objetivos\") {\n        $enObjectivesParagraph = $paragraphs.Item($i + 2)\n        break\n    }\n}\nif ($enObjectivesParagraph -eq $null) {\n    throw \"Could not locate empty English objectives paragraph\"\n}\n$enObjectivesParagraph.Range.InsertBefore('Develop in students the skills and abilities necessary to apply scientific, technological and engineering knowledge in the conception, design, installation, optimization, supervision and critical evaluation of the operation of bioprocesses, with an emphasis on: 1) Types of bioreactors; 2) Bioreactor operating modes and 3) Bioreactor analysis.')\n\n# 2) Straightforward text replacements via Find/Replace.\n$replacements = @(\n    @{ Old = 'Ativa\u00e7\u00e3o: 01/01/2018'; New = 'Ativa\u00e7\u00e3o: 01/01/2025' },\n    @{ Old = 'Capacitar o aluno para identificar e aplicar os conceitos de Reatores Bioqu\u00edmicos em diferentes bioprocessos (fermentativos e enzim\u00e1ticos). Especificamente, capacitar o aluno para definir os tipos de biorreatores; para definir as diferentes formas de condu\u00e7\u00e3o de um processo fermentativo empregando biorreatores, para realizar o equacionamento matem\u00e1tico do crescimento microbiano e da forma\u00e7\u00e3o de produtos de interesse de um bioprocesso empregando biorreatores e para aplica\u00e7\u00e3o dos conceitos gerais dos reatores enzim\u00e1ticos operados em diferentes fases.'; New = 'Desenvolver nos discentes as compet\u00eancias e habilidades necess\u00e1rias para a aplica\u00e7\u00e3o de conhecimentos cient\u00edficos, tecnol\u00f3gicos e de engenharia na concep\u00e7\u00e3o, projeto, instala\u00e7\u00e3o, otimiza\u00e7\u00e3o, supervis\u00e3o e avalia\u00e7\u00e3o cr\u00edtica da opera\u00e7\u00e3o de bioprocessos, com \u00eanfase em: 1) Tipos de biorreatores; 2) Formas de opera\u00e7\u00e3o dos biorreatores e 3) An\u00e1lise de biorreatores' },\n    @{ Old = 'Introdu\u00e7\u00e3o a biorreatores; processo descont\u00ednuo; processo cont\u00ednuo; processo descont\u00ednuo alimentado e reatores enzim\u00e1ticos.'; New = 'Ser\u00e3o apresentados os principais tipos de biorreatores associados as suas aplica\u00e7\u00f5es para diferentes bioprocessos; definidas as formas de opera\u00e7\u00e3o do biorreator e analisadas as diferentes formas de opera\u00e7\u00e3o de biorreatores com base nos balan\u00e7os materiais dos componentes do sistema.' },\n    @{ Old = 'Introduction to bioreactors; batch process; continuous process; fed batch process and enzymatic reactors.'; New = 'The main types of bioreactors associated with their applications for different bioprocesses will be presented. Also, the main forms of bioreactor operation will be defined and further analyzed based on the material balances of the system components.' },\n    @{ Old = '1. Introdu\u00e7\u00e3o a biorreatores: apresenta\u00e7\u00e3o e classifica\u00e7\u00e3o de reatores bioqu\u00edmicos; intera\u00e7\u00e3o microorganismos/meios para estabelecimento de condi\u00e7\u00f5es para c\u00e1lculos de biorreatores.2. Processo descont\u00ednuo: caracter\u00edsticas gerais do processo cont\u00ednuo; balan\u00e7o material para c\u00e9lula, substrato e produto, em um \u00fanico est\u00e1gio com e sem reciclo de c\u00e9lulas; aplica\u00e7\u00e3o do processo cont\u00ednuo (exemplos).3. Processo cont\u00ednuo: caracter\u00edsticas gerais do processo cont\u00ednuo; balan\u00e7o material para c\u00e9lula, substrato e produto, em um \u00fanico est\u00e1gio com e sem reciclo de c\u00e9lulas; aplica\u00e7\u00e3o do processo cont\u00ednuo (exemplos).4. Processo descont\u00ednuo alimentado: caracter\u00edsticas gerais do processo descont\u00ednuo alimentado; balan\u00e7o material para c\u00e9lula e substrato, com volume vari\u00e1vel, empregando vaz\u00e3o constante de alimenta\u00e7\u00e3o; considera\u00e7\u00f5es sobre forma\u00e7\u00e3o de produtos no processo descont\u00ednuo alimentado; aplica\u00e7\u00e3o do processo descont\u00ednuo alimentado (exemplos).5. Reatores enzim\u00e1ticos: caracter\u00edsticas gerais dos reatores enzim\u00e1ticos; aplica\u00e7\u00e3o de processos enzim\u00e1ticos (exemplos).'; New = '1) Defini\u00e7\u00e3o e classifica\u00e7\u00e3o de biorreatores; 2) Biorreatores para cultivos submersos (agitados mecanicamente, agitados pneumaticamente, leito fixo, leito fluidizado, outros tipos; 3) Biorreatores para cultivos em estado s\u00f3lido (est\u00e1ticos e agitados); 4) Principais formas de opera\u00e7\u00e3o de biorreatores (descont\u00ednua, cont\u00ednua e descont\u00ednua-alimentada) e 4) An\u00e1lise de biorreatores (balan\u00e7os materiais para c\u00e9lulas, substrato-limitante e produtos metab\u00f3licos) nas diferentes formas de opera\u00e7\u00e3o.' },\n    @{ Old = '1.Introduction to bioreactors: Description and classification of biochemical reactors; bioreactors configuration; microorganisms/medium interaction; bioreactor operation modes.2.Batch process: general characteristics of the batch process; material balance for cell, substrate and product; application of batch process (examples).3.Continuous process: general characteristics, operation modes (single or multiple stages with or without cells recycle); material balance for cell and substrate; formation of products in the continuous systems, application of continuous process (examples).4.Fed batch process: general characteristics; material balance for cell and substrate (equations for fed-batch operation at variable and fixed volume); considerations about formation of products in the fed batch process; application of fed batch process (examples).5.Enzymatic reactors: general characteristics of the enzymatic reactors; application of enzymatic processes (examples).'; New = '1) Definition and classification of bioreactors; 2) Bioreactors for submerged cultures (mechanically agitated, pneumatically agitated, fixed bed, fluidized bed, other types; 3) Bioreactors for solid state cultures (static and agitated); 4) Main forms of bioreactor operation (batch, continuous and fed-batch) and 4) Analysis of bioreactors (material balances for cells, limiting substrate and metabolic products) in different forms of operation.' },\n    @{ Old = 'Os alunos ser\u00e3o avaliados formalmente por duas provas te\u00f3ricas. A pondera\u00e7\u00e3o das notas ser\u00e1 de 50% para cada avalia\u00e7\u00e3o, ou seja: M\u00e9dia do per\u00edodo letivo normal = (P1 + P2 )/2'; New = 'Os alunos ser\u00e3o avaliados formalmente por duas provas te\u00f3ricas. A pondera\u00e7\u00e3o das notas ser\u00e1 de 50% para cada avalia\u00e7\u00e3o, ou seja: M\u00e9dia do per\u00edodo letivo normal = (P1 + P2)/2.' },\n    @{ Old = '1)  Aiba, S., Humphrey, A.E., Millis, N.F. Biochemical Engineering - 2\u00aa Edi\u00e7\u00e3o- 1973.2) Asenjo A., Merchuk, J.C. Bioreactor System Design-1995.3) Stanbury, D. and Whitaker, A. Principles af Fermentation Technology-1986.4) Lima, U.A., Aquarone, E., Borzani, W. Biotecnologia Industrial. Fundamentos Vol. 1, Engenharia Bioqu\u00edmica Vol.2, Processos Fermentativos Vol.3. Ed.  Edgard Blucher, S\u00e3o Paulo, 2001.'; New = 'ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioqu\u00edmica. 2\u00aa Edi\u00e7\u00e3o. S\u00e3o Paulo: Blucher, 2021. ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914PESSOA JR, A; VITOLO, M; LONG, P.F.(editors).  Pharmaceutical Biotechnology: A Focus on Industrial Application. CRC Press, 1st Edition.2022.' },\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $ok) {\n        throw \"Replacement not applied: $($r.Old.Substring(0, [Math]::Min(40, $r.Old.Length)))\"\n    }\n}\n\nWrite-Output \"done\"\n"}
